$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 616.9167
$ws.Cells.Item(6, 9).Value = 267
$ws.Cells.Item(6, 11).Value = 801
$ws.Cells.Item(6, 13).Value = -689

$ws.Cells.Item(21, 8).Value = 10602.333
$ws.Cells.Item(21, 10).Value = 10602.333
$ws.Cells.Item(21, 12).Value = 10602.333
$ws.Cells.Item(21, 14).Value = -11538.333

$ws.Cells.Item(23, 8).Value = 10602.333
$ws.Cells.Item(23, 10).Value = 10602.333
$ws.Cells.Item(23, 12).Value = 10602.333
$ws.Cells.Item(23, 14).Value = -11070.333

$ws.Cells.Item(34, 8).Value = 18314.666
$ws.Cells.Item(34, 9).Value = 2472
$ws.Cells.Item(34, 11).Value = 2472
$ws.Cells.Item(34, 13).Value = -2269

$ws.Cells.Item(36, 8).Value = 18314.666
$ws.Cells.Item(36, 9).Value = 2472
$ws.Cells.Item(36, 11).Value = 2472
$ws.Cells.Item(36, 13).Value = -1757

$ws.Cells.Item(137, 8).Value = 32259244
$ws.Cells.Item(137, 9).Value = 45455372
$ws.Cells.Item(137, 10).Value = 2039.4445
$ws.Cells.Item(137, 11).Value = 136366116
$ws.Cells.Item(137, 12).Value = 6118.333500000001
$ws.Cells.Item(137, 13).Value = -136363566
$ws.Cells.Item(137, 14).Value = -11218.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 18948.75
$ws.Cells.Item(63, 9).Value = 22865
$ws.Cells.Item(63, 11).Value = 22865
$ws.Cells.Item(63, 13).Value = -22179

$ws.Cells.Item(66, 8).Value = 18948.75
$ws.Cells.Item(66, 9).Value = 22865
$ws.Cells.Item(66, 11).Value = 114325
$ws.Cells.Item(66, 13).Value = -110893

$ws.Cells.Item(109, 8).Value = 23500
$ws.Cells.Item(109, 10).Value = 23500
$ws.Cells.Item(109, 12).Value = 23500
$ws.Cells.Item(109, 14).Value = -26274

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1284.7778
$ws.Cells.Item(94, 9).Value = 945.375
$ws.Cells.Item(94, 11).Value = 945.375
$ws.Cells.Item(94, 13).Value = -494.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5599.5557
$ws.Cells.Item(31, 9).Value = 2035.8948
$ws.Cells.Item(31, 10).Value = 9582.471
$ws.Cells.Item(31, 11).Value = 2035.8948
$ws.Cells.Item(31, 12).Value = 9582.471
$ws.Cells.Item(31, 13).Value = -1740.8948
$ws.Cells.Item(31, 14).Value = -10172.471

$ws.Cells.Item(34, 8).Value = 5599.5557
$ws.Cells.Item(34, 9).Value = 2035.8948
$ws.Cells.Item(34, 10).Value = 9582.471
$ws.Cells.Item(34, 11).Value = 2035.8948
$ws.Cells.Item(34, 12).Value = 9582.471
$ws.Cells.Item(34, 13).Value = -1833.8948
$ws.Cells.Item(34, 14).Value = -9986.471

$ws.Cells.Item(62, 8).Value = 17483.467
$ws.Cells.Item(62, 9).Value = 26283.334
$ws.Cells.Item(62, 10).Value = 4283.6665
$ws.Cells.Item(62, 11).Value = 26283.334
$ws.Cells.Item(62, 12).Value = 4283.6665
$ws.Cells.Item(62, 13).Value = -25659.334
$ws.Cells.Item(62, 14).Value = -5531.6665

$ws.Cells.Item(65, 8).Value = 17483.467
$ws.Cells.Item(65, 9).Value = 26283.334
$ws.Cells.Item(65, 10).Value = 4283.6665
$ws.Cells.Item(65, 11).Value = 131416.67
$ws.Cells.Item(65, 12).Value = 21418.3325
$ws.Cells.Item(65, 13).Value = -128296.67
$ws.Cells.Item(65, 14).Value = -27658.3325

$ws.Cells.Item(99, 8).Value = 4814824.5
$ws.Cells.Item(99, 9).Value = 5215793
$ws.Cells.Item(99, 10).Value = 3200
$ws.Cells.Item(99, 11).Value = 5215793
$ws.Cells.Item(99, 12).Value = 3200
$ws.Cells.Item(99, 13).Value = -5214295
$ws.Cells.Item(99, 14).Value = -6196

$ws.Cells.Item(122, 8).Value = 1525.2
$ws.Cells.Item(122, 9).Value = 1281.5
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 3844.5
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = -1394.5
$ws.Cells.Item(122, 14).Value = -12400

$ws.Cells.Item(126, 8).Value = 4814824.5
$ws.Cells.Item(126, 9).Value = 5215793
$ws.Cells.Item(126, 10).Value = 3200
$ws.Cells.Item(126, 11).Value = 15647379
$ws.Cells.Item(126, 12).Value = 9600
$ws.Cells.Item(126, 13).Value = -15644909
$ws.Cells.Item(126, 14).Value = -14540

$ws.Cells.Item(135, 8).Value = 44993.332
$ws.Cells.Item(135, 10).Value = 44993.332
$ws.Cells.Item(135, 12).Value = 44993.332
$ws.Cells.Item(135, 14).Value = -55133.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 58.04348
$ws.Cells.Item(12, 10).Value = 39.555557
$ws.Cells.Item(12, 12).Value = 118.666671
$ws.Cells.Item(12, 14).Value = -464.666671

$ws.Cells.Item(98, 8).Value = 434
$ws.Cells.Item(98, 10).Value = 302
$ws.Cells.Item(98, 12).Value = 906
$ws.Cells.Item(98, 14).Value = -3902

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6704.5386
$ws.Cells.Item(70, 9).Value = 6817.6665
$ws.Cells.Item(70, 10).Value = 6450
$ws.Cells.Item(70, 11).Value = 6817.6665
$ws.Cells.Item(70, 12).Value = 6450
$ws.Cells.Item(70, 13).Value = -6547.6665
$ws.Cells.Item(70, 14).Value = -6990

$ws.Cells.Item(73, 8).Value = 6704.5386
$ws.Cells.Item(73, 9).Value = 6817.6665
$ws.Cells.Item(73, 10).Value = 6450
$ws.Cells.Item(73, 11).Value = 6817.6665
$ws.Cells.Item(73, 12).Value = 6450
$ws.Cells.Item(73, 13).Value = -5881.6665
$ws.Cells.Item(73, 14).Value = -8322

$ws.Cells.Item(110, 8).Value = 56635.332
$ws.Cells.Item(110, 10).Value = 56635.332
$ws.Cells.Item(110, 12).Value = 56635.332
$ws.Cells.Item(110, 14).Value = -64815.332

$ws.Cells.Item(126, 8).Value = 3430.3447
$ws.Cells.Item(126, 9).Value = 2064.4443
$ws.Cells.Item(126, 10).Value = 4045
$ws.Cells.Item(126, 11).Value = 6193.3329
$ws.Cells.Item(126, 12).Value = 12135
$ws.Cells.Item(126, 13).Value = -3723.3329
$ws.Cells.Item(126, 14).Value = -17075

$ws.Cells.Item(132, 8).Value = 2287.6667
$ws.Cells.Item(132, 9).Value = 2397.55
$ws.Cells.Item(132, 10).Value = 2067.9
$ws.Cells.Item(132, 11).Value = 7192.650000000001
$ws.Cells.Item(132, 12).Value = 6203.700000000001
$ws.Cells.Item(132, 13).Value = -4662.650000000001
$ws.Cells.Item(132, 14).Value = -11263.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 8655.77
$ws.Cells.Item(22, 9).Value = 856
$ws.Cells.Item(22, 10).Value = 21135.4
$ws.Cells.Item(22, 11).Value = 856
$ws.Cells.Item(22, 12).Value = 21135.4
$ws.Cells.Item(22, 13).Value = -561
$ws.Cells.Item(22, 14).Value = -21725.4

$ws.Cells.Item(27, 8).Value = 8655.77
$ws.Cells.Item(27, 9).Value = 856
$ws.Cells.Item(27, 10).Value = 21135.4
$ws.Cells.Item(27, 11).Value = 856
$ws.Cells.Item(27, 12).Value = 21135.4
$ws.Cells.Item(27, 13).Value = -749
$ws.Cells.Item(27, 14).Value = -21349.4

$ws.Cells.Item(93, 8).Value = 1431.3846
$ws.Cells.Item(93, 9).Value = 1161.3636
$ws.Cells.Item(93, 10).Value = 2916.5
$ws.Cells.Item(93, 11).Value = 1161.3636
$ws.Cells.Item(93, 12).Value = 2916.5
$ws.Cells.Item(93, 13).Value = 86.63640000000009
$ws.Cells.Item(93, 14).Value = -5412.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 955812.6
$ws.Cells.Item(81, 9).Value = 2002086.6
$ws.Cells.Item(81, 10).Value = 4654.4546
$ws.Cells.Item(81, 11).Value = 4004173.2
$ws.Cells.Item(81, 12).Value = 9308.9092
$ws.Cells.Item(81, 13).Value = -4003112.2
$ws.Cells.Item(81, 14).Value = -11430.9092

$ws.Cells.Item(84, 8).Value = 955812.6
$ws.Cells.Item(84, 9).Value = 2002086.6
$ws.Cells.Item(84, 10).Value = 4654.4546
$ws.Cells.Item(84, 11).Value = 20020866
$ws.Cells.Item(84, 12).Value = 46544.546
$ws.Cells.Item(84, 13).Value = -20015562
$ws.Cells.Item(84, 14).Value = -57152.546

$ws.Cells.Item(126, 8).Value = 92118.63
$ws.Cells.Item(126, 9).Value = 125938.125
$ws.Cells.Item(126, 10).Value = 1933.3334
$ws.Cells.Item(126, 11).Value = 377814.375
$ws.Cells.Item(126, 12).Value = 5800.0002
$ws.Cells.Item(126, 13).Value = -375344.375
$ws.Cells.Item(126, 14).Value = -10740.0002
